# Commit: "I0 and IF added" -- append two new data columns (I: "I0",
# J: "IF") to the stats sheet, with header cells formatted like the
# existing headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting used by B1:H1 (bold font, thin box
# border, centered horizontal alignment, top vertical alignment).
$headerSrc = $ws.Range("H1")
$headerDst = $ws.Range("I1:J1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2..73: fill columns I and J ---
$I = @(6,8,5,6,7,6,7,7,8,5,8,7,9,8,5,5,6,9,10,6,7,6,6,5,9,6,7,6,6,6,6,5,6,6,8,6,8,6,9,5,6,7,5,6,9,7,8,7,9,7,6,7,8,6,4,9,6,5,5,9,7,9,6,8,8,6,5,5,6,6,4,3)
$J = @(6,8,6,6,7,6,7,7,8,6,9,7,9,8,5,6,7,9,10,6,7,6,7,5,10,7,7,6,7,7,6,5,6,6,8,6,8,6,9,6,6,7,5,6,9,7,8,7,9,7,7,8,8,6,5,9,6,5,7,9,8,9,6,8,8,6,6,5,6,6,4,3)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
